# Updated cryptos list on Thu Jan  4 19:55:03 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.100.81'
$ws.Range('E2').Value = '  +2.95%  '
$ws.Range('D3').Value = '2.277.68'
$ws.Range('E3').Value = '  +3.01%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = '''318.16'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').Value = '''106.12'
$ws.Range('E6').Value = '  +7.97%  '
$ws.Range('E7').Value = '  +1.38%  '
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('E9').Value = '  +2.43%  '
$ws.Range('D10').Value = '''38.98'
$ws.Range('D11').Value = '''0.0841'
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').Value = '''7.91'
$ws.Range('E12').Value = '  +1.86%  '
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('D14').Value = '2.625.80'
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('D15').Value = '''0.881'
$ws.Range('E15').Value = '  +2.53%  '
$ws.Range('D16').Value = '''14.63'
$ws.Range('E16').Value = '  +3.48%  '
$ws.Range('D17').Value = '2.277.75'
$ws.Range('E17').Value = '  +3.09%  '
$ws.Range('D18').Value = '44.041.28'
$ws.Range('E18').Value = '  +3.01%  '
$ws.Range('E19').Value = '  -4.86%  '
$ws.Range('D20').Value = '0.0₃0999'
$ws.Range('E20').Value = '  +4.48%  '
$ws.Range('D21').Value = '''6.55'
$ws.Range('E21').Value = '  +2.90%  '
$ws.Range('D22').Value = '''66.23'
$ws.Range('E22').Value = '  +1.72%  '
$ws.Range('D23').Value = '''3.21'
$ws.Range('E23').Value = '  +1.59%  '
$ws.Range('D24').Value = '''237.82'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').Value = '''2.20'
$ws.Range('E25').Value = '  +4.64%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '''10.27'
$ws.Range('E27').Value = '  +2.30%  '
$ws.Range('D28').Value = '''38.96'
$ws.Range('E28').Value = '  +15.44%  '
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').Value = '''6.53'
$ws.Range('E30').Value = '  +4.17%  '
$ws.Range('D31').Value = '''163.42'
$ws.Range('E31').Value = '  +4.51%  '
$ws.Range('D32').Value = '''20.61'
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('D33').Value = '''0.0882'
$ws.Range('E33').Value = '  +1.63%  '
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '''2.08'
$ws.Range('E35').Value = '  +4.36%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').Value = '''3.27'
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('D37').Value = '''0.115'
$ws.Range('E37').Value = '  +12.49%  '
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('D39').Value = '''3.96'
$ws.Range('E39').Value = '  +8.19%  '
$ws.Range('D40').Value = '''4.49'
$ws.Range('E40').Value = '  +1.31%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '''15.52'
$ws.Range('E41').Value = '  +28.18%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0326'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('D44').Value = '1.766.88'
$ws.Range('E44').Value = '  -6.62%  '
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('D46').Value = '''85.83'
$ws.Range('E46').Value = '  -3.95%  '
$ws.Range('D47').Value = '''5.39'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = '''75.31'
$ws.Range('E48').Value = '  +0.78%  '
$ws.Range('D49').Value = '''8.80'
$ws.Range('E49').Value = '  +2.90%  '
$ws.Range('D50').Value = '''59.67'
$ws.Range('E50').Value = '  -1.01%  '
$ws.Range('D51').Value = '''104.50'
$ws.Range('E51').Value = '  +3.34%  '
